# Xhosa translation of the Child Safety interview information & consent form.
#
# Word's Find.Execute(..., Replace:=wdReplaceOne/All) runs text through the
# smart-quotes AutoCorrect pass (straight ' -> curly ') and the engine only
# re-emits xml:space="preserve" on a <w:t> when the *final* text itself has
# leading/trailing whitespace. Assigning directly to a Range's .Text avoids
# the AutoCorrect rewrite, and when that Range is an entire single-run
# paragraph's Range (i.e. it still includes the paragraph mark) the
# xml:space="preserve" hint on that run is preserved too. So: use Find only
# to *locate* each paragraph, then mutate via Range.Text.

$d = $word.ActiveDocument

function Set-ParagraphText($needle, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $needle"
    }
    $para = $rng.Paragraphs(1)
    $para.Range.Text = $newText
}

# 1) Title
Set-ParagraphText `
    "Appendix 19: SWIFT Child Safety Module Interview: Information Sheet and Consent Form" `
    "ISihlomelo 19: Udliwano-ndlebe lweModyuli yoKhuseleko loMntwana ye-SWIFT: Iphepha loLwazi kunye neFomu yeMvume"

# 2) "Your experience with the ParentText programme..." paragraph
Set-ParagraphText `
    "Your experience with the ParentText programme is vital to our study." `
    "Amava akho nenkqubo ye-ParentText abalulekile kuphononongo lwethu. Singathanda ukuva ngamava wakho emodyuli ngoku 'Gcina Umntwana Wakho Ekhuselekile'. Oludliwano-ndlebe luyinxalenye yophononongo olwenziwa ngabaphandi abaphuma kwiDyunivesithi yaseKapa kunye neyoMzantsi Afrika kunye neDyunivesithi yaseOxford eUnited Kingdom. "

# 3) "Before you decide if you'd like to be interviewed..." paragraph
Set-ParagraphText `
    "Before you decide if you" `
    "Ngaphambili kokuba ugqibe ekubeni ungathanda na ukuba nodliwano-ndlebe, kubalulekile ukuba wazi kutheni sisenza oluphando nje kwaye ukuthatha inxalenye kungaquka ntoni. Lonke ulwazi ozakuludinga luchaziwe ngezantsi kodwa ukuba unayo nayiphi na imibuzo malunga nokuthatha inxaxheba okanye ngophando lwethu, ndicela u-imeyilele iqela lophononongo ku swift@globalparenting.org okanye uthumele umyalezo kuthi ku WhatsApp at +27 XX XXX XXXX. Silapha ukuzokunceda wena!"

# 4) " What will my interview look like and what is expected of me?" heading
Set-ParagraphText `
    "What will my interview look like and what is expected of me?" `
    " Udliwano-ndlebe lwam luzakujongeka njani kwaye kulindeleke ntoni kum?"

# 5) "We would like to have a telephonic conversation..." paragraph
Set-ParagraphText `
    "We would like to have a telephonic conversation with you which will last a maximum of 45 minutes." `
    "Singathanda ukuba nencoko ngomxeba nawe ozakuthatha imizuzu engamashumi amane anesihlanu ubude. Omnye wabaphandi bethu uzakutsalela umnxeba athethe nawe ngexesha elikulungeleyo wena. Akukho zimpendulu zilungileyo okanye ezingalunganga, sifuna nje ukuva amava kunye nemibono yakho nge chatbot. Please make sure that when we call, that you only let the interview start when you are in a private space where you feel comfortable to talk without being overheard or interrupted. If while you are being interviewed, you are interrupted, please ask the researcher to pause until you feel safe to continue talking."
